# Update the "Fitness" column (C) values in the log sheet.
# Rows 2-128   -> 7534
# Rows 129-143 -> 7345
# Rows 144-170 -> 7343
# Rows 171-222 -> 7293

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 128; $r++) {
    $ws.Cells.Item($r, 3).Value = 7534
}

for ($r = 129; $r -le 143; $r++) {
    $ws.Cells.Item($r, 3).Value = 7345
}

for ($r = 144; $r -le 170; $r++) {
    $ws.Cells.Item($r, 3).Value = 7343
}

for ($r = 171; $r -le 222; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
